$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (CloneScene entry) - entire row shifts everything up
$ws.Rows("2:2").Delete()

# Delete row 3 (was RebellerNoob entry, now at row 3 after first delete) - shifts up again
$ws.Rows("3:3").Delete()

# Update RelivePos for new row2 (was PioneerNoob/villageScene row)
$ws.Range("E2").Value = "20,0,60"

# Update ID for new row3 (was the blank-FilePath Demo1 row)
$ws.Range("B3").Value = "2"

# Match the author's final selection/active cell in the sheet view
[void]$ws.Range("F5").Select()
